$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4142.857
$ws.Range("H43").Value = 4659.8
$ws.Range("J43").Value = 4999.75
$ws.Range("L43").Value = 4999.75
$ws.Range("N43").Value = -5137.75
$ws.Range("H86").Value = 10009.529
$ws.Range("I86").Value = 8941.333000000001
$ws.Range("K86").Value = 8941.333000000001
$ws.Range("M86").Value = -7818.333000000001
$ws.Range("H89").Value = 10009.529
$ws.Range("I89").Value = 8941.333000000001
$ws.Range("K89").Value = 44706.665
$ws.Range("M89").Value = -39090.665
$ws.Range("H113").Value = 4014.0715
$ws.Range("J113").Value = 4285.5713
$ws.Range("L113").Value = 4285.5713
$ws.Range("N113").Value = -10793.5713
$ws.Range("H132").Value = 9366.5
$ws.Range("I132").Value = 9990.454
$ws.Range("K132").Value = 29971.362
$ws.Range("M132").Value = -27441.362
$ws.Range("H133").Value = 99990
$ws.Range("J133").Value = 99990
$ws.Range("L133").Value = 99990
$ws.Range("N133").Value = -110110
$ws.Range("H138").Value = 2112.0278
$ws.Range("I138").Value = 1858.5807
$ws.Range("K138").Value = 5575.742099999999
$ws.Range("M138").Value = -435.7420999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 1800
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -1671
$ws.Range("N33").Value = -1658
$ws.Range("H45").Value = 2643.6667
$ws.Range("I45").Value = 1776.625
$ws.Range("K45").Value = 1776.625
$ws.Range("M45").Value = -1399.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37207
$ws.Range("I20").Value = 54028.668
$ws.Range("J20").Value = 1881.5
$ws.Range("K20").Value = 54028.668
$ws.Range("L20").Value = 1881.5
$ws.Range("M20").Value = -53781.668
$ws.Range("N20").Value = -2375.5
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9664
$ws.Range("N33").ClearContents()
$ws.Range("H99").Value = 10740.363
$ws.Range("I99").Value = 10740.363
$ws.Range("K99").Value = 10740.363
$ws.Range("M99").Value = -9242.362999999999
$ws.Range("H107").Value = 1496.875
$ws.Range("I107").Value = 1269.35
$ws.Range("J107").Value = 1876.0834
$ws.Range("K107").Value = 1269.35
$ws.Range("L107").Value = 1876.0834
$ws.Range("M107").Value = 650.6500000000001
$ws.Range("N107").Value = -5716.0834
$ws.Range("H134").Value = 1010.4583
$ws.Range("I134").Value = 888.1905
$ws.Range("K134").Value = 2664.5715
$ws.Range("M134").Value = -129.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3181.3333
$ws.Range("I31").Value = 2955.6667
$ws.Range("K31").Value = 2955.6667
$ws.Range("M31").Value = -2660.6667
$ws.Range("H34").Value = 3181.3333
$ws.Range("I34").Value = 2955.6667
$ws.Range("K34").Value = 2955.6667
$ws.Range("M34").Value = -2753.6667
$ws.Range("H99").Value = 13099.55
$ws.Range("I99").Value = 22933.2
$ws.Range("K99").Value = 22933.2
$ws.Range("M99").Value = -21435.2
$ws.Range("H122").Value = 21197.77
$ws.Range("I122").Value = 1745.8572
$ws.Range("K122").Value = 5237.571599999999
$ws.Range("M122").Value = -2787.571599999999
$ws.Range("H126").Value = 13099.55
$ws.Range("I126").Value = 22933.2
$ws.Range("K126").Value = 68799.60000000001
$ws.Range("M126").Value = -66329.60000000001
$ws.Range("H132").Value = 4293.647
$ws.Range("I132").Value = 5741.0967
$ws.Range("J132").Value = 2050.1
$ws.Range("K132").Value = 17223.2901
$ws.Range("L132").Value = 6150.299999999999
$ws.Range("M132").Value = -14693.2901
$ws.Range("N132").Value = -11210.3
$ws.Range("H134").Value = 1597.862
$ws.Range("I134").Value = 1301.6666
$ws.Range("K134").Value = 3904.9998
$ws.Range("M134").Value = -1369.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 550
$ws.Range("I86").Value = 700
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 2100
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -914
$ws.Range("N86").Value = -3572
$ws.Range("H89").Value = 550
$ws.Range("I89").Value = 700
$ws.Range("J89").Value = 400
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 3600
$ws.Range("M89").Value = -372
$ws.Range("N89").Value = -15456

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 13500
$ws.Range("J26").Value = 13500
$ws.Range("L26").Value = 13500
$ws.Range("N26").Value = -14060
$ws.Range("H50").Value = 13500
$ws.Range("J50").Value = 13500
$ws.Range("L50").Value = 13500
$ws.Range("N50").Value = -14496
$ws.Range("H80").Value = 2468.6667
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 2453
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2453
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -4449
$ws.Range("H83").Value = 2468.6667
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 2453
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 12265
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -22249
$ws.Range("H102").Value = 2936.5715
$ws.Range("I102").Value = 1134.3334
$ws.Range("K102").Value = 1134.3334
$ws.Range("M102").Value = 487.6666
$ws.Range("H122").Value = 3126.353
$ws.Range("I122").Value = 2324.75
$ws.Range("J122").Value = 3373
$ws.Range("K122").Value = 6974.25
$ws.Range("L122").Value = 10119
$ws.Range("M122").Value = -4524.25
$ws.Range("N122").Value = -15019
$ws.Range("H132").Value = 6781.2324
$ws.Range("I132").Value = 5240.1
$ws.Range("J132").Value = 10337.692
$ws.Range("K132").Value = 15720.3
$ws.Range("L132").Value = 31013.076
$ws.Range("M132").Value = -13190.3
$ws.Range("N132").Value = -36073.076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 4333.3335
$ws.Range("I32").Value = 4333.3335
$ws.Range("K32").Value = 4333.3335
$ws.Range("M32").Value = -4016.3335
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 25000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 25000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -24710
$ws.Range("N33").ClearContents()
$ws.Range("H55").Value = 1308.6786
$ws.Range("I55").Value = 1470.2727
$ws.Range("J55").Value = 1204.1177
$ws.Range("K55").Value = 1470.2727
$ws.Range("L55").Value = 1204.1177
$ws.Range("M55").Value = -1297.2727
$ws.Range("N55").Value = -1550.1177
$ws.Range("H58").Value = 6348.25
$ws.Range("I58").Value = 2796.5
$ws.Range("J58").Value = 9900
$ws.Range("K58").Value = 2796.5
$ws.Range("L58").Value = 9900
$ws.Range("M58").Value = -2536.5
$ws.Range("N58").Value = -10420
$ws.Range("H68").Value = 4362.593
$ws.Range("I68").Value = 3400.0557
$ws.Range("K68").Value = 3400.0557
$ws.Range("M68").Value = -2651.0557
$ws.Range("H71").Value = 4362.593
$ws.Range("I71").Value = 3400.0557
$ws.Range("K71").Value = 17000.2785
$ws.Range("M71").Value = -13256.2785
$ws.Range("H93").Value = 5251.8184
$ws.Range("I93").Value = 3835.1
$ws.Range("J93").Value = 19419
$ws.Range("K93").Value = 3835.1
$ws.Range("L93").Value = 19419
$ws.Range("M93").Value = -2587.1
$ws.Range("N93").Value = -21915
$ws.Range("H122").Value = 6317.75
$ws.Range("I122").Value = 3298.5
$ws.Range("K122").Value = 9895.5
$ws.Range("M122").Value = -7445.5
$ws.Range("H132").Value = 3291.7693
$ws.Range("I132").Value = 3285.4285
$ws.Range("J132").Value = 3299.1667
$ws.Range("K132").Value = 9856.2855
$ws.Range("L132").Value = 9897.500100000001
$ws.Range("M132").Value = -7326.2855
$ws.Range("N132").Value = -14957.5001
$ws.Range("H136").Value = 4666.3335
$ws.Range("I136").Value = 4666.3335
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13999.0005
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11449.0005
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7199.273
$ws.Range("I62").Value = 6812.7144
$ws.Range("J62").Value = 7875.75
$ws.Range("K62").Value = 6812.7144
$ws.Range("L62").Value = 7875.75
$ws.Range("M62").Value = -6188.7144
$ws.Range("N62").Value = -9123.75
$ws.Range("H65").Value = 7199.273
$ws.Range("I65").Value = 6812.7144
$ws.Range("J65").Value = 7875.75
$ws.Range("K65").Value = 34063.572
$ws.Range("L65").Value = 39378.75
$ws.Range("M65").Value = -30943.572
$ws.Range("N65").Value = -45618.75
$ws.Range("H122").Value = 55839.617
$ws.Range("I122").Value = 1793
$ws.Range("K122").Value = 5379
$ws.Range("M122").Value = -2929
$ws.Range("H130").Value = 39600
$ws.Range("J130").Value = 39600
$ws.Range("L130").Value = 39600
$ws.Range("N130").Value = -49640
$ws.Range("H132").Value = 2191.037
$ws.Range("I132").Value = 2083
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6249
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3719
$ws.Range("N132").Value = -20060
